# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G17").Value = 3000000
$ws.Range("G19").Value = 1600000
$ws.Range("G20").Value = 2000000
